# Auto-generated edit script applying the scheduled-runner data refresh
# to the Cuchulainn_Profits (FFXIV market/leve) tracking workbook.
# For each sheet, cells in columns H-N are overwritten with refreshed
# market data. Cells that no longer have a value (diff removes them)
# are cleared so no stale cached value remains.

$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H8").Value = 25.461538
$ws.Range("I8").Value = 13.888889
$ws.Range("K8").Value = 41.666667
$ws.Range("M8").Value = 97.333333
$ws.Range("H69").Value = 3507.5
$ws.Range("I69").Value = 5000
$ws.Range("K69").Value = 15000
$ws.Range("M69").Value = -14126
$ws.Range("H70").Value = 3700.25
$ws.Range("I70").Value = 2150
$ws.Range("J70").Value = 5250.5
$ws.Range("K70").Value = 6450
$ws.Range("L70").Value = 15751.5
$ws.Range("M70").Value = -6180
$ws.Range("N70").Value = -16291.5
$ws.Range("H72").Value = 3507.5
$ws.Range("I72").Value = 5000
$ws.Range("K72").Value = 45000
$ws.Range("M72").Value = -40632
$ws.Range("H73").Value = 3700.25
$ws.Range("I73").Value = 2150
$ws.Range("J73").Value = 5250.5
$ws.Range("K73").Value = 6450
$ws.Range("L73").Value = 15751.5
$ws.Range("M73").Value = -5514
$ws.Range("N73").Value = -17623.5
$ws.Range("J113").Value = 3500
$ws.Range("L113").Value = 3500
$ws.Range("N113").Value = -10008
$ws.Range("H129").Value = 5348.5
$ws.Range("I129").Value = 1197
$ws.Range("K129").Value = 3591
$ws.Range("M129").Value = 1409
$ws.Range("H135").Value = 1790
$ws.Range("I135").Value = 1790
$ws.Range("K135").Value = 16110
$ws.Range("M135").Value = -13575

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("H16").Value = 185.42857
$ws.Range("I16").Value = 26
$ws.Range("J16").Value = 398
$ws.Range("K16").Value = 26
$ws.Range("L16").Value = 398
$ws.Range("M16").Value = 261
$ws.Range("N16").Value = -972
$ws.Range("H26").Value = 4333.3335
$ws.Range("I26").Value = 4333.3335
$ws.Range("K26").Value = 4333.3335
$ws.Range("M26").Value = -4003.3335
$ws.Range("H39").Value = 7000
$ws.Range("I39").Value = 10000
$ws.Range("J39").Value = 4000
$ws.Range("K39").Value = 10000
$ws.Range("L39").Value = 4000
$ws.Range("M39").Value = -9480
$ws.Range("N39").Value = -5040
$ws.Range("H45").Value = 6930
$ws.Range("I45").Value = 4408.75
$ws.Range("K45").Value = 4408.75
$ws.Range("M45").Value = -4031.75
$ws.Range("H97").Value = 1254
$ws.Range("I97").Value = 1498
$ws.Range("J97").Value = 1010
$ws.Range("K97").Value = 1498
$ws.Range("L97").Value = 1010
$ws.Range("M97").Value = -1002
$ws.Range("N97").Value = -2002
$ws.Range("H122").Value = 1005
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1005
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 3015
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -7915
$ws.Range("H132").Value = 1566.6666
$ws.Range("I132").Value = 1566.6666
$ws.Range("K132").Value = 4699.9998
$ws.Range("M132").Value = -2169.9998

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H20").Value = 3120.1428
$ws.Range("I20").Value = 2998.5
$ws.Range("J20").Value = 3168.8
$ws.Range("K20").Value = 2998.5
$ws.Range("L20").Value = 3168.8
$ws.Range("M20").Value = -2751.5
$ws.Range("N20").Value = -3662.8
$ws.Range("H99").Value = 1100
$ws.Range("I99").Value = 1100
$ws.Range("K99").Value = 1100
$ws.Range("M99").Value = 398
$ws.Range("H105").Value = 1570.5
$ws.Range("I105").Value = 1570.5
$ws.Range("K105").Value = 1570.5
$ws.Range("M105").Value = 176.5
$ws.Range("H134").Value = 3334.5557
$ws.Range("I134").Value = 2876.375
$ws.Range("K134").Value = 8629.125
$ws.Range("M134").Value = -6094.125

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H2").Value = 340.4
$ws.Range("I2").Value = 69
$ws.Range("J2").Value = 747.5
$ws.Range("K2").Value = 69
$ws.Range("L2").Value = 747.5
$ws.Range("M2").Value = 44
$ws.Range("N2").Value = -973.5
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()
$ws.Range("H13").Value = 975
$ws.Range("J13").Value = 975
$ws.Range("L13").Value = 975
$ws.Range("N13").Value = -1253
$ws.Range("H35").Value = 1842.2858
$ws.Range("J35").Value = 2748
$ws.Range("L35").Value = 2748
$ws.Range("N35").Value = -3336
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H122").Value = 1070
$ws.Range("I122").Value = 600
$ws.Range("K122").Value = 1800
$ws.Range("M122").Value = 650

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H132").Value = 1387.375
$ws.Range("I132").Value = 566.3333
$ws.Range("J132").Value = 1880
$ws.Range("K132").Value = 5096.9997
$ws.Range("L132").Value = 16920
$ws.Range("M132").Value = -2566.9997
$ws.Range("N132").Value = -21980

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H70").Value = 5749.8335
$ws.Range("I70").Value = 5374.75
$ws.Range("K70").Value = 5374.75
$ws.Range("M70").Value = -5104.75
$ws.Range("H73").Value = 5749.8335
$ws.Range("I73").Value = 5374.75
$ws.Range("K73").Value = 5374.75
$ws.Range("M73").Value = -4438.75
$ws.Range("H97").Value = 600
$ws.Range("J97").Value = 600
$ws.Range("L97").Value = 600
$ws.Range("N97").Value = -1592
$ws.Range("H122").Value = 1670
$ws.Range("I122").Value = 1670
$ws.Range("K122").Value = 5010
$ws.Range("M122").Value = -2560

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H12").Value = 10700
$ws.Range("J12").Value = 1400
$ws.Range("L12").Value = 1400
$ws.Range("N12").Value = -1740
$ws.Range("H68").Value = 2800
$ws.Range("I68").Value = 2400
$ws.Range("K68").Value = 2400
$ws.Range("M68").Value = -1651
$ws.Range("H71").Value = 2800
$ws.Range("I71").Value = 2400
$ws.Range("K71").Value = 12000
$ws.Range("M71").Value = -8256
$ws.Range("H122").Value = 3411.5
$ws.Range("I122").Value = 3073
$ws.Range("J122").Value = 3750
$ws.Range("K122").Value = 9219
$ws.Range("L122").Value = 11250
$ws.Range("M122").Value = -6769
$ws.Range("N122").Value = -16150

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H3").Value = 1000
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 1000
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 1000
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -1228
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("H34").Value = 1000
$ws.Range("J34").Value = 1000
$ws.Range("L34").Value = 1000
$ws.Range("N34").Value = -1406
$ws.Range("H62").Value = 69999.336
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 69999.336
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 69999.336
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -71247.336
$ws.Range("H65").Value = 69999.336
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 69999.336
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 349996.68
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -356236.68
$ws.Range("H122").Value = 2799.8333
$ws.Range("I122").Value = 2799.8333
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8399.499899999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5949.499899999999
$ws.Range("N122").ClearContents()
$ws.Range("H136").Value = 7957.5713
$ws.Range("J136").Value = 6874.5
$ws.Range("L136").Value = 20623.5
$ws.Range("N136").Value = -25723.5

